# Saldo_guide.xlsx update
# - Reference date for every data row (column G, "Dt. Referencia") moves
#   forward one day: 2024-07-03 (serial 45476) -> 2024-07-04 (serial 45477).
# - A handful of "Saldo Previsto" / "Vl. Total" amounts (columns E and H)
#   were corrected for specific accounts.
# - The worksheet (tab) name is refreshed to match the new extraction
#   timestamp embedded in its name.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the "Dt. Referencia" column (G) for every data row (2-275) forward
# by one day in a single broadcast assignment.
$ws.Range("G2:G275").Value = 45477

# Corrected amounts: same new value written to both "Saldo Previsto" (E)
# and "Vl. Total" (H) for the affected accounts.
$corrections = @{
    52  = 319.92
    55  = 119.55
    57  = 513.86
    110 = 6326.64
    112 = 936.44
    113 = 95.9
    165 = 788.61
    255 = 458.32
}

foreach ($row in $corrections.Keys) {
    $value = $corrections[$row]
    $ws.Range("E$row").Value = $value
    $ws.Range("H$row").Value = $value
}

# Refresh the worksheet name to reflect the new export timestamp.
$ws.Name = "IClientBalance-20240704-104725-"
